# Regenerate merged AHB files
# Remove the "ÄNDERUNG" marker (shared string) from column L for a set of
# rows, resetting those cells back to the plain centered/grey style (style
# index 4 in the original workbook) instead of the highlighted "ÄNDERUNG"
# style (style index 7).
#
# Also fix two "group header" rows (212 and 218) whose formatting had
# drifted to the generic body-row style (index 5 / 5) instead of the
# header-row style used by every other first-row-of-a-group (index 2 for
# most columns, index 3 for column B) - matching the established pattern
# seen e.g. in rows 2, 9, 12, 16, ...

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Source cell that already carries the desired plain "s=4" style used for
# column L on normal (non "ÄNDERUNG") rows.
$formatSourceL = $ws.Range("L2")

# Contiguous row ranges (in column L) that only need the "ÄNDERUNG" value
# cleared and their style reset to the plain style.
$lRanges = @(
    "L80",
    "L98:L106",
    "L108",
    "L111:L114",
    "L117:L136",
    "L139:L140",
    "L143",
    "L149:L151",
    "L158:L159",
    "L165:L169",
    "L175",
    "L187:L189",
    "L196",
    "L223:L224"
)

$formatSourceL.Copy()
foreach ($rng in $lRanges) {
    $ws.Range($rng).PasteSpecial($xlPasteFormats)
}
foreach ($rng in $lRanges) {
    $ws.Range($rng).Value = ""
}

# Rows 212 and 218: these are the first row of a new group (like rows 2, 9,
# 12, ...) and should use the "header" row formatting instead of the
# generic body-row formatting. Copy formats column-by-column from row 2,
# which already has the correct target formatting pattern (style 2 for
# most columns, style 3 for column B), then clear column L's value same as
# the other rows above.
$headerRows = @(212, 218)

foreach ($r in $headerRows) {
    $ws.Range("A2").Copy()
    $ws.Range("A$r").PasteSpecial($xlPasteFormats)

    $ws.Range("B2").Copy()
    $ws.Range("B$r").PasteSpecial($xlPasteFormats)

    $ws.Range("C2:V2").Copy()
    $ws.Range("C$r`:V$r").PasteSpecial($xlPasteFormats)

    $ws.Range("L$r").Value = ""
}

$excel.CutCopyMode = 0
